$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $val) {
    $ws.Range($cellRef).Value = "'" + $val
}

$ws.Range("D2").Value = "59.775.10"
$ws.Range("E2").Value = "  -1.88%  "
$ws.Range("D3").Value = "2.604.46"
$ws.Range("E3").Value = "  +0.36%  "
$ws.Range("E4").Value = "  -0.07%  "
Set-TextCell "D5" "514.20"
$ws.Range("E5").Value = "  -1.78%  "
$ws.Range("E6").Value = "  -5.49%  "
$ws.Range("E7").Value = "  +0.09%  "
Set-TextCell "D8" "0.563"
$ws.Range("E8").Value = "  -4.92%  "
$ws.Range("D9").Value = "2.606.86"
$ws.Range("E9").Value = "  +0.19%  "
Set-TextCell "D10" "6.25"
$ws.Range("E10").Value = "  -6.91%  "
$ws.Range("E11").Value = "  -2.03%  "
Set-TextCell "D12" "0.336"
$ws.Range("E12").Value = "  -3.48%  "
$ws.Range("E13").Value = "  -0.88%  "
$ws.Range("D14").Value = "3.060.20"
$ws.Range("E14").Value = "  +0.29%  "
$ws.Range("D15").Value = "59.762.93"
$ws.Range("E15").Value = "  -1.93%  "
$ws.Range("E16").Value = "  -3.48%  "
$ws.Range("E17").Value = "  -3.54%  "
$ws.Range("D18").Value = "2.601.90"
$ws.Range("E18").Value = "  +0.20%  "
Set-TextCell "D19" "4.58"
$ws.Range("E19").Value = "  -3.55%  "
Set-TextCell "D20" "338.39"
$ws.Range("E20").Value = "  -4.18%  "
Set-TextCell "D21" "10.25"
$ws.Range("E21").Value = "  -3.21%  "
Set-TextCell "D22" "6.03"
$ws.Range("E22").Value = "  -3.30%  "
$ws.Range("E23").Value = "  -0.09%  "
Set-TextCell "D24" "60.67"
$ws.Range("E24").Value = "  -0.74%  "
Set-TextCell "D25" "0.413"
$ws.Range("E25").Value = "  -3.25%  "
Set-TextCell "D26" "0.999"
$ws.Range("E26").Value = "  +0.15%  "
$ws.Range("E27").Value = "  -5.27%  "
$ws.Range("E28").Value = "  -6.57%  "
Set-TextCell "D29" "6.94"
$ws.Range("E29").Value = "  -5.95%  "
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("E31").Value = "  -2.35%  "
Set-TextCell "D32" "5.91"
$ws.Range("E32").Value = "  -6.31%  "
Set-TextCell "D35" "3.89"
$ws.Range("E35").Value = "  -7.23%  "
Set-TextCell "D36" "0.895"
$ws.Range("E36").Value = "  -4.33%  "
$ws.Range("E37").Value = "  -6.88%  "
Set-TextCell "D38" "36.61"
$ws.Range("E38").Value = "  +0.36%  "
Set-TextCell "D39" "0.839"
$ws.Range("E39").Value = "  -1.37%  "
$ws.Range("E40").Value = "  -6.61%  "
Set-TextCell "D41" "3.57"
$ws.Range("E41").Value = "  -6.06%  "
Set-TextCell "D42" "284.40"
$ws.Range("E42").Value = "  -1.17%  "
Set-TextCell "D43" "0.620"
$ws.Range("E43").Value = "  -0.43%  "
Set-TextCell "D44" "0.999"
$ws.Range("E44").Value = "  +0.16%  "
$ws.Range("E45").Value = "  -2.20%  "
Set-TextCell "D46" "0.0539"
$ws.Range("E46").Value = "  -3.91%  "
$ws.Range("E47").Value = "  -3.97%  "
Set-TextCell "D48" "10.37"
$ws.Range("E48").Value = "  +0.50%  "
Set-TextCell "D49" "0.0230"
$ws.Range("E49").Value = "  -3.48%  "

# Rows 33/34 and 50/51: coin identity swap
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextCell "D33" "18.72"
$ws.Range("E33").Value = "  -3.37%  "

$ws.Range("B34").Value = "Monero"
$ws.Range("C34").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell "D34" "150.17"
$ws.Range("E34").Value = "  +1.27%  "

$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell "D50" "4.56"
$ws.Range("E50").Value = "  -6.57%  "

$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "1.931.35"
$ws.Range("E51").Value = "  -1.09%  "
